# Insert 9 new rows of historical data (2019-11-18 .. 2019-11-28) for stock
# 5271 / PECCA right after the existing 2019-11-15 row, pushing all the
# following rows down by 9 (old row 875 "2019-11-29" becomes row 884, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything currently at row 875 downward by 9 rows, leaving 9
# blank rows at 875:883 for the new data.
$ws.Rows("875:883").Insert()

# Columns B (date) and C (id) hold values that look like dates/numbers
# ("2019-11-18", "5271") but must stay plain text, matching the rest of
# the sheet. Mark them as Text *before* writing the values so Excel does
# not auto-convert them to a date serial / number.
$ws.Range("B875:C883").NumberFormat = "@"

# row 875 : 2019-11-18
$ws.Range("A875").Value = 1574035200
$ws.Range("B875").Value = "2019-11-18"
$ws.Range("C875").Value = "5271"
$ws.Range("D875").Value = "PECCA"
$ws.Range("E875").Value = 1.22
$ws.Range("F875").Value = 1.22
$ws.Range("G875").Value = 1.18
$ws.Range("H875").Value = 1.19
$ws.Range("I875").Value = 155500

# row 876 : 2019-11-19
$ws.Range("A876").Value = 1574121600
$ws.Range("B876").Value = "2019-11-19"
$ws.Range("C876").Value = "5271"
$ws.Range("D876").Value = "PECCA"
$ws.Range("E876").Value = 1.18
$ws.Range("F876").Value = 1.22
$ws.Range("G876").Value = 1.18
$ws.Range("H876").Value = 1.2
$ws.Range("I876").Value = 358700

# row 877 : 2019-11-20
$ws.Range("A877").Value = 1574208000
$ws.Range("B877").Value = "2019-11-20"
$ws.Range("C877").Value = "5271"
$ws.Range("D877").Value = "PECCA"
$ws.Range("E877").Value = 1.21
$ws.Range("F877").Value = 1.21
$ws.Range("G877").Value = 1.17
$ws.Range("H877").Value = 1.19
$ws.Range("I877").Value = 409600

# row 878 : 2019-11-21
$ws.Range("A878").Value = 1574294400
$ws.Range("B878").Value = "2019-11-21"
$ws.Range("C878").Value = "5271"
$ws.Range("D878").Value = "PECCA"
$ws.Range("E878").Value = 1.19
$ws.Range("F878").Value = 1.19
$ws.Range("G878").Value = 1.17
$ws.Range("H878").Value = 1.18
$ws.Range("I878").Value = 329800

# row 879 : 2019-11-22
$ws.Range("A879").Value = 1574380800
$ws.Range("B879").Value = "2019-11-22"
$ws.Range("C879").Value = "5271"
$ws.Range("D879").Value = "PECCA"
$ws.Range("E879").Value = 1.19
$ws.Range("F879").Value = 1.22
$ws.Range("G879").Value = 1.19
$ws.Range("H879").Value = 1.22
$ws.Range("I879").Value = 278300

# row 880 : 2019-11-25
$ws.Range("A880").Value = 1574640000
$ws.Range("B880").Value = "2019-11-25"
$ws.Range("C880").Value = "5271"
$ws.Range("D880").Value = "PECCA"
$ws.Range("E880").Value = 1.22
$ws.Range("F880").Value = 1.22
$ws.Range("G880").Value = 1.19
$ws.Range("H880").Value = 1.2
$ws.Range("I880").Value = 560000

# row 881 : 2019-11-26
$ws.Range("A881").Value = 1574726400
$ws.Range("B881").Value = "2019-11-26"
$ws.Range("C881").Value = "5271"
$ws.Range("D881").Value = "PECCA"
$ws.Range("E881").Value = 1.2
$ws.Range("F881").Value = 1.22
$ws.Range("G881").Value = 1.19
$ws.Range("H881").Value = 1.22
$ws.Range("I881").Value = 118900

# row 882 : 2019-11-27
$ws.Range("A882").Value = 1574812800
$ws.Range("B882").Value = "2019-11-27"
$ws.Range("C882").Value = "5271"
$ws.Range("D882").Value = "PECCA"
$ws.Range("E882").Value = 1.22
$ws.Range("F882").Value = 1.27
$ws.Range("G882").Value = 1.21
$ws.Range("H882").Value = 1.24
$ws.Range("I882").Value = 1971000

# row 883 : 2019-11-28
$ws.Range("A883").Value = 1574899200
$ws.Range("B883").Value = "2019-11-28"
$ws.Range("C883").Value = "5271"
$ws.Range("D883").Value = "PECCA"
$ws.Range("E883").Value = 1.24
$ws.Range("F883").Value = 1.25
$ws.Range("G883").Value = 1.23
$ws.Range("H883").Value = 1.25
$ws.Range("I883").Value = 153400
